$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Update "想去人数" (F column) counters that bumped up ---
    $ws.Range("F2").Value  = 1074
    $ws.Range("F5").Value  = 4617
    $ws.Range("F6").Value  = 27
    $ws.Range("F7").Value  = 387
    $ws.Range("F8").Value  = 1367
    $ws.Range("F9").Value  = 906
    $ws.Range("F11").Value = 1045
    $ws.Range("F13").Value = 572

    # --- Insert a new event row at row 15, pushing the existing       ---
    # --- "南昌·CM01动漫游戏博览会" (old row 15) and                   ---
    # --- "鹰潭·原×铁×崩only" (old row 16) rows down by one.           ---
    $ws.Rows.Item(15).Insert()

    # Copy formatting (bold/border/centering) from the row above so the
    # newly-inserted row's index cell matches the rest of column A.
    $ws.Range("A14").Copy()
    $ws.Range("A15").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # New row 15: 南昌·运动番only春季集训
    $ws.Range("A15").Value = 14
    # Force text so "2024.03.23" isn't auto-converted into a date serial,
    # then restore the plain (unstyled) format used by the other B cells.
    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = "2024.03.23"
    $ws.Range("B14").Copy()
    $ws.Range("B15").PasteSpecial(-4122)
    $excel.CutCopyMode = 0
    $ws.Range("C15").Value = "南昌·运动番only春季集训"
    $ws.Range("D15").Value = "创新三路777号 南昌小飞侠章鱼文化体育公园"
    $ws.Range("E15").Value = "2024.03.23 10:00-03.24 17:00"
    $ws.Range("F15").Value = 3
    $ws.Range("G15").Value = 58
    $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=81950"
    $ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"

    # Row 16 (was row 15 before the insert): update the "想去人数" count
    $ws.Range("A16").Value = 15
    $ws.Range("F16").Value = 262

    # Row 17 (was row 16 before the insert): renumber the index column
    $ws.Range("A17").Value = 16
}
